# Auto-generated edit script: updates cryptos list values per commit
# "Updated cryptos list on Mon Mar 18 02:05:17 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "67.673.51"
$ws.Cells.Item(2, 5).Value = "  +2.28%  "
$ws.Cells.Item(3, 4).Value = "3.602.33"
$ws.Cells.Item(3, 5).Value = "  +1.28%  "
$ws.Cells.Item(4, 5).Value = "  +0.32%  "
$ws.Cells.Item(5, 4).Value = "'199.13"
$ws.Cells.Item(5, 5).Value = "  +6.77%  "
$ws.Cells.Item(6, 4).Value = "'558.58"
$ws.Cells.Item(6, 5).Value = "  -4.28%  "
$ws.Cells.Item(7, 4).Value = "3.600.88"
$ws.Cells.Item(7, 5).Value = "  +1.46%  "
$ws.Cells.Item(8, 4).Value = "'0.612"
$ws.Cells.Item(8, 5).Value = "  +0.22%  "
$ws.Cells.Item(9, 5).Value = "  +0.17%  "
$ws.Cells.Item(10, 4).Value = "'0.669"
$ws.Cells.Item(10, 5).Value = "  +0.69%  "
$ws.Cells.Item(11, 4).Value = "'57.48"
$ws.Cells.Item(11, 5).Value = "  +7.73%  "
$ws.Cells.Item(12, 5).Value = "  +4.37%  "
$ws.Cells.Item(13, 4).Value = "'0.0000288"
$ws.Cells.Item(13, 5).Value = "  +13.36%  "
$ws.Cells.Item(14, 4).Value = "'9.96"
$ws.Cells.Item(14, 5).Value = "  +2.79%  "
$ws.Cells.Item(15, 4).Value = "4.186.91"
$ws.Cells.Item(15, 5).Value = "  +1.85%  "
$ws.Cells.Item(16, 4).Value = "3.601.98"
$ws.Cells.Item(16, 5).Value = "  +1.51%  "
$ws.Cells.Item(17, 5).Value = "  +0.58%  "
$ws.Cells.Item(18, 4).Value = "'18.80"
$ws.Cells.Item(18, 5).Value = "  +3.28%  "
$ws.Cells.Item(19, 4).Value = "67.589.93"
$ws.Cells.Item(19, 5).Value = "  +2.55%  "
$ws.Cells.Item(20, 5).Value = "  +1.03%  "
$ws.Cells.Item(21, 4).Value = "'1.07"
$ws.Cells.Item(21, 5).Value = "  +2.48%  "
$ws.Cells.Item(22, 4).Value = "'399.13"
$ws.Cells.Item(22, 5).Value = "  +1.54%  "
$ws.Cells.Item(23, 4).Value = "'12.99"
$ws.Cells.Item(23, 5).Value = "  +22.34%  "
$ws.Cells.Item(24, 4).Value = "'4.12"
$ws.Cells.Item(24, 5).Value = "  -4.79%  "
$ws.Cells.Item(25, 4).Value = "'84.89"
$ws.Cells.Item(25, 5).Value = "  +0.05%  "
$ws.Cells.Item(26, 4).Value = "'2.94"
$ws.Cells.Item(26, 5).Value = "  +2.70%  "
$ws.Cells.Item(27, 4).Value = "'12.42"
$ws.Cells.Item(27, 5).Value = "  +0.24%  "
$ws.Cells.Item(28, 4).Value = "'6.10"
$ws.Cells.Item(28, 5).Value = "  +1.25%  "
$ws.Cells.Item(29, 4).Value = "'3.83"
$ws.Cells.Item(29, 5).Value = "  +8.66%  "
$ws.Cells.Item(30, 4).Value = "'8.40"
$ws.Cells.Item(30, 5).Value = "  +23.70%  "
$ws.Cells.Item(31, 4).Value = "'9.24"
$ws.Cells.Item(31, 5).Value = "  +3.98%  "
$ws.Cells.Item(32, 4).Value = "'31.41"
$ws.Cells.Item(32, 5).Value = "  +2.34%  "
$ws.Cells.Item(33, 4).Value = "'672.33"
$ws.Cells.Item(33, 5).Value = "  +9.73%  "
$ws.Cells.Item(34, 5).Value = "  +0.89%  "
$ws.Cells.Item(35, 4).Value = "'0.113"
$ws.Cells.Item(35, 5).Value = "  +2.37%  "
$ws.Cells.Item(36, 4).Value = "'63.46"
$ws.Cells.Item(36, 5).Value = "  +0.67%  "
$ws.Cells.Item(37, 4).Value = "'42.35"
$ws.Cells.Item(37, 5).Value = "  +3.04%  "
$ws.Cells.Item(38, 4).Value = "'0.430"
$ws.Cells.Item(38, 5).Value = "  +15.72%  "
$ws.Cells.Item(39, 5).Value = "  -0.09%  "
$ws.Cells.Item(40, 4).Value = "0.0₃0772"
$ws.Cells.Item(40, 5).Value = "  +4.16%  "
$ws.Cells.Item(41, 4).Value = "'3.19"
$ws.Cells.Item(41, 5).Value = "  +15.36%  "
$ws.Cells.Item(42, 2).Value = "Maker"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(42, 4).Value = "3.227.90"
$ws.Cells.Item(42, 5).Value = "  +10.31%  "
$ws.Cells.Item(43, 2).Value = "Kaspa"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(43, 4).Value = "'0.137"
$ws.Cells.Item(43, 5).Value = "  +4.78%  "
$ws.Cells.Item(44, 2).Value = "Fetch.AI"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(44, 4).Value = "'2.82"
$ws.Cells.Item(44, 5).Value = "  +15.63%  "
$ws.Cells.Item(45, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(45, 4).Value = "'0.999"
$ws.Cells.Item(45, 5).Value = "  +0.13%  "
$ws.Cells.Item(46, 2).Value = "dogwifhat"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(46, 4).Value = "'3.00"
$ws.Cells.Item(46, 5).Value = "  +29.04%  "
$ws.Cells.Item(47, 4).Value = "'0.0415"
$ws.Cells.Item(47, 5).Value = "  +2.95%  "
$ws.Cells.Item(48, 4).Value = "'2.70"
$ws.Cells.Item(48, 5).Value = "  +10.54%  "
$ws.Cells.Item(49, 2).Value = "Stellar"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(49, 4).Value = "'0.130"
$ws.Cells.Item(49, 5).Value = "  +0.39%  "
$ws.Cells.Item(50, 2).Value = "THORChain"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Cells.Item(50, 4).Value = "'8.71"
$ws.Cells.Item(50, 5).Value = "  +2.38%  "
$ws.Cells.Item(51, 2).Value = "ApeXProtocol"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(51, 4).Value = "'3.08"
$ws.Cells.Item(51, 5).Value = "  -0.10%  "
